$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4027.411126602437
$ws.Range("C3").Value = 4027.411126602437
$ws.Range("C4").Value = 3950.581415222695
$ws.Range("C5").Value = 3950.581415222695
$ws.Range("C6").Value = 3947.81345803351
$ws.Range("C7").Value = 3947.81345803351
$ws.Range("C8").Value = 3947.81345803351
$ws.Range("C9").Value = 3947.81345803351
$ws.Range("C10").Value = 3904.353809646612
$ws.Range("C11").Value = 3805.897927000234
$ws.Range("C12").Value = 3805.897927000234
